$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# Prepare row 46 styles by copying single cells from row 45 (keeps
# the copy rectangle 1 column wide so no stray empty cells appear
# in columns that should stay absent, e.g. J46/L46).
# -----------------------------------------------------------------
$ws.Range("A45").Copy()
$ws.Range("A46").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("I45").Copy()
$ws.Range("I46").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A46").Value = 42862.298773148148

# -----------------------------------------------------------------
# Row 46 text cells that reuse already existing shared strings.
# -----------------------------------------------------------------
$ws.Range("B46").Value = "            Buy"
$ws.Range("C46").Value = "        XRP"

# -----------------------------------------------------------------
# New shared strings have to be introduced in this exact order so
# they land on the expected indices (matches how the workbook was
# originally authored): fee(E) -> pair code(G) -> price(D) -> amount(F)
# -----------------------------------------------------------------
$ws.Range("E46").Value = "         0.096USDT"

$ws.Range("G45").Value = " XRP/USDT0000004"
$ws.Range("G46").Value = " XRP/USDT0000004"

# D46 holds a numeric-looking string, so force text formatting,
# assign the value, then restore the original (wrap-text / general)
# number format copied from D45 so the resulting cell style matches.
$ws.Range("D45").Copy()
$ws.Range("D46").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "             0.10560001`n"
$ws.Range("D45").Copy()
$ws.Range("D46").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("F46").Value = "         410 XRP"

$ws.Range("H46").Value = "IN PROGRESS"
$ws.Range("K46").Value = "     "

# -----------------------------------------------------------------
# Row 45: the old "IN PROGRESS" trade is finalized as CANCEL, now
# referencing the new trading pair label, with a finalized timestamp.
# -----------------------------------------------------------------
$ws.Range("H45").Value = "CANCEL"
$ws.Range("I45").Value = 42862.298773148148
$ws.Range("K45").Value = "     "

# -----------------------------------------------------------------
# Keep row heights consistent with the rest of the sheet (assigning
# the multi-line D46 text can otherwise auto-grow the row).
# -----------------------------------------------------------------
$ws.Rows.Item(45).RowHeight = 14.25
$ws.Rows.Item(46).RowHeight = 14.25

# -----------------------------------------------------------------
# Selection / scroll position, matching the saved view state.
# -----------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 29
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D47").Select()
